$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was added to the "Perejil" (Feria Lagunitas de
# Puerto Montt) data set. It belongs chronologically right after the row
# that is currently row 185, so insert a fresh row at position 186 and
# push the existing rows 186:293 down to 187:294.
$ws.Rows.Item(186).Insert()

# Populate the newly inserted row 186 with the new record's data.
$ws.Cells.Item(186, 1).Value = 4
$ws.Cells.Item(186, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(186, 3).Value = "Los Lagos"
$ws.Cells.Item(186, 4).Value = 44830
$ws.Cells.Item(186, 5).Value = 10
$ws.Cells.Item(186, 6).Value = 100112044
$ws.Cells.Item(186, 7).Value = "Perejil"
$ws.Cells.Item(186, 8).Value = "Sin especificar"
$ws.Cells.Item(186, 9).Value = "Primera"
$ws.Cells.Item(186, 10).Value = 80
$ws.Cells.Item(186, 11).Value = 6000
$ws.Cells.Item(186, 12).Value = 6000
$ws.Cells.Item(186, 13).Value = 6000
$ws.Cells.Item(186, 14).Value = "`$/docena de atados (3 kilos)"
$ws.Cells.Item(186, 15).Value = "Región Metropolitana"
$ws.Cells.Item(186, 16).Value = 2000
$ws.Cells.Item(186, 17).Value = 3
$ws.Cells.Item(186, 18).Value = "Hortaliza"

# Match the date formatting style used by the rest of the column D cells.
$ws.Cells.Item(186, 4).NumberFormat = $ws.Cells.Item(187, 4).NumberFormat
